$d = $word.ActiveDocument

# --- Merge the proofErr-split runs back into single runs (no visible text change) ---
# Paragraph 1: "**Document Name* ... Voltage_current  "
$rng1 = $d.Paragraphs.Item(1).Range
$rng1.Find.Execute(" Voltage_current  ", $false, $false, $false, $false, $false, $true, 1, $false, " Voltage_current  ", 2) | Out-Null

# Paragraph 2: "**Gpt canvas* ..."
$rng2 = $d.Paragraphs.Item(2).Range
$rng2.Find.Execute("**Gpt canvas*", $false, $false, $false, $false, $false, $true, 1, $false, "**Gpt canvas*", 2) | Out-Null

# --- Remove the "speaker notes" related paragraphs (feature removed) ---
# Work from the highest paragraph index down so earlier indices stay valid:
#   7 -> trailing empty paragraph
#   6 -> "**Ppt with speaker notes** : N"
#   3 -> "**Gpt canvas without speaker notes* : N"
$d.Paragraphs.Item(7).Range.Delete() | Out-Null
$d.Paragraphs.Item(6).Range.Delete() | Out-Null
$d.Paragraphs.Item(3).Range.Delete() | Out-Null
